$d = $word.ActiveDocument

# Die MwSt-Angabe wurde von 19% auf 16% geändert.
# Alle Vorkommen von "zuzüglich 19 % MwSt" in den Tabellenzellen ersetzen.
$d.Content.Find.Execute(
    "zuzüglich 19 % MwSt", $true, $false, $false, $false, $false,
    $true, 1, $false, "zuzüglich 16 % MwSt", 2
)

# Das automatisch gepflegte "_GoBack"-Lesezeichen (letzte Bearbeitungsposition)
# am Ende des Dokuments entfernen, falls vorhanden.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
